$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update timestamps
$wsOverview.Range("G2").Value = "2016-08-30 02:41:35"
$wsDeDe.Range("H2").Value = "2016-08-30 02:41:35"
$wsZhCn.Range("H2").Value = "2016-08-30 02:41:31"

# Update column widths
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
